# "Major update. Schedule, New Parameters." commit
#
# The functional change in Tests/Test Parameters.xlsx is the new default
# "From" day parameter used by the schedule/termination logic: cell C2
# (named range FROM) moves from 11 to 31. C3 (named range TO) is untouched.
#
# (The rest of the published diff - new xr/xr6/xr10/xr2 revision
# namespaces, fileVersion/rupBuild bump, absPath, workbookView size,
# defaultRowHeight/dyDescent and the sub-pixel column width - are just
# artifacts of the authoring machine re-saving the file with a newer
# Excel build; they are not reachable/settable through the Excel object
# model, so we leave them alone rather than risk corrupting the sheet.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "From" day parameter.
$ws.Range("C2").Value = 31

# Row heights for the bordered rows grew slightly (13.5 -> 13.8pt) in the
# resave; this is exactly reproducible, so apply it for fidelity.
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 13.8
